$wb = $excel.ActiveWorkbook

# --- Training Dashboard sheet ---
$ws1 = $wb.Worksheets.Item("Training Dashboard")

$ws1.Range("H3").Value = 638
$ws1.Range("I3").Formula = "'16-Sep-2025"

$ws1.Range("H4").Value = 378
$ws1.Range("I4").Formula = "'16-Sep-2025"

$ws1.Range("H5").Value = 435
$ws1.Range("I5").Formula = "'16-Sep-2025"

$ws1.Range("H6").Value = 423
$ws1.Range("I6").Formula = "'16-Sep-2025"

$ws1.Range("H7").Value = 170
$ws1.Range("I7").Formula = "'16-Sep-2025"

$ws1.Range("H8").Value = 254
$ws1.Range("I8").Formula = "'16-Sep-2025"

$ws1.Range("H9").Value = 331
$ws1.Range("I9").Formula = "'16-Sep-2025"

# --- Exam Dashboard sheet ---
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

$ws2.Columns.Item(5).ColumnWidth = 15

$ws2.Range("E3").Value = "date is valid"
$ws2.Range("E4").Value = "date is valid"
